$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.891.93'
$ws.Range("E2").Value = '  -0.86%  '

$ws.Range("D3").Value = '1.630.40'
$ws.Range("E3").Value = '  -2.26%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '''206.23'
$ws.Range("E5").Value = '  -1.82%  '

$ws.Range("D6").Value = '''0.5120'
$ws.Range("E6").Value = '  -1.68%  '

$ws.Range("D7").Value = '''1.005'
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").Value = '''0.2539'
$ws.Range("E8").Value = '  -3.69%  '

$ws.Range("D9").Value = '''0.06154'
$ws.Range("E9").Value = '  -1.16%  '

$ws.Range("D10").Value = '''20.29'
$ws.Range("E10").Value = '  -4.26%  '

$ws.Range("D11").Value = '''0.07554'
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").Value = '1.631.90'
$ws.Range("E12").Value = '  -2.76%  '

$ws.Range("D13").Value = '''4.331'
$ws.Range("E13").Value = '  -2.03%  '

$ws.Range("D14").Value = '1.852.75'
$ws.Range("E14").Value = '  -2.28%  '

$ws.Range("D15").Value = '''0.5337'
$ws.Range("E15").Value = '  -4.56%  '

$ws.Range("D16").Value = '0.0₅7927'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").Value = '''65.04'
$ws.Range("E17").Value = '  -1.97%  '

$ws.Range("D18").Value = '25.941.32'
$ws.Range("E18").Value = '  -0.92%  '

$ws.Range("D19").Value = '''1.005'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = '''4.595'
$ws.Range("E20").Value = '  -3.89%  '

$ws.Range("D21").Value = '''184.68'
$ws.Range("E21").Value = '  -1.05%  '

$ws.Range("D22").Value = '''9.933'
$ws.Range("E22").Value = '  -4.08%  '

$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("D24").Value = '''6.036'
$ws.Range("E24").Value = '  -2.07%  '

$ws.Range("D25").Value = '''146.69'
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("D26").Value = '''0.1189'
$ws.Range("E26").Value = '  -4.35%  '

$ws.Range("D27").Value = '''7.260'
$ws.Range("E27").Value = '  -4.05%  '

$ws.Range("D28").Value = '''15.37'
$ws.Range("E28").Value = '  -3.32%  '

$ws.Range("D29").Value = '''1.353'
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").Value = '''0.05981'
$ws.Range("E30").Value = '  -4.07%  '

$ws.Range("D31").Value = '''1.240'
$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").Value = '''3.379'
$ws.Range("E32").Value = '  -2.85%  '

$ws.Range("D33").Value = '''3.341'
$ws.Range("E33").Value = '  -2.40%  '

$ws.Range("D34").Value = '''1.603'
$ws.Range("E34").Value = '  -1.41%  '

$ws.Range("D35").Value = '''0.9609'
$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("D36").Value = '''2.384'
$ws.Range("E36").Value = '  -0.82%  '

$ws.Range("D37").Value = '''2.715'
$ws.Range("E37").Value = '  +0.43%  '

$ws.Range("D38").Value = '''0.5777'
$ws.Range("E38").Value = '  -4.33%  '

$ws.Range("D39").Value = '''0.01575'
$ws.Range("E39").Value = '  -1.82%  '

$ws.Range("D40").Value = '1.071.19'
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("D41").Value = '''5.764'
$ws.Range("E41").Value = '  -5.89%  '

$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").Value = '''0.8385'
$ws.Range("E43").Value = '  -3.12%  '

$ws.Range("D44").Value = '''99.53'
$ws.Range("E44").Value = '  +0.37%  '

$ws.Range("D45").Value = '1.785.23'
$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -2.75%  '

$ws.Range("D47").Value = '''0.9983'
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("D48").Value = '''53.98'
$ws.Range("E48").Value = '  -3.54%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05201'
$ws.Range("E49").Value = '  -0.98%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.896'
$ws.Range("E50").Value = '  -0.64%  '

$ws.Range("D51").Value = '''0.4230'
$ws.Range("E51").Value = '  -0.51%  '
